$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H12").Value = 267.5
$ws.Range("I12").Value = 269.44446
$ws.Range("K12").Value = 269.44446
$ws.Range("M12").Value = -99.44445999999999
$ws.Range("H28").Value = 948.25
$ws.Range("J28").Value = 798.5
$ws.Range("L28").Value = 798.5
$ws.Range("N28").Value = -1768.5
$ws.Range("H106").Value = 5614.4346
$ws.Range("I106").Value = 5581.85
$ws.Range("K106").Value = 5581.85
$ws.Range("M106").Value = -4950.85
$ws.Range("H137").Value = 3804.6382
$ws.Range("I137").Value = 3830.6
$ws.Range("J137").Value = 3728.9167
$ws.Range("K137").Value = 11491.8
$ws.Range("L137").Value = 11186.7501
$ws.Range("M137").Value = -8941.799999999999
$ws.Range("N137").Value = -16286.7501

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H21").Value = 618.3333
$ws.Range("I21").Value = 618.3333
$ws.Range("K21").Value = 618.3333
$ws.Range("M21").Value = -244.3333
$ws.Range("H32").Value = 6505.073
$ws.Range("I32").Value = 8655.241
$ws.Range("K32").Value = 8655.241
$ws.Range("M32").Value = -8368.241
$ws.Range("H122").Value = 6001.2104
$ws.Range("I122").Value = 5918.6665
$ws.Range("K122").Value = 17755.9995
$ws.Range("M122").Value = -15305.9995
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H128").Value = 70429
$ws.Range("J128").Value = 70429
$ws.Range("L128").Value = 70429
$ws.Range("N128").Value = -80389
$ws.Range("H132").Value = 1939.5294
$ws.Range("I132").Value = 1749.5
$ws.Range("J132").Value = 4980
$ws.Range("K132").Value = 5248.5
$ws.Range("L132").Value = 14940
$ws.Range("M132").Value = -2718.5
$ws.Range("N132").Value = -20000
$ws.Range("H135").Value = 40000
$ws.Range("J135").Value = 40000
$ws.Range("L135").Value = 40000
$ws.Range("N135").Value = -50140
$ws.Range("H137").Value = 49166.668
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200
$ws.Range("H139").Value = 45000
$ws.Range("I139").Value = 20000
$ws.Range("K139").Value = 20000
$ws.Range("M139").Value = -14860

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H22").Value = 766.6667
$ws.Range("I22").Value = 501
$ws.Range("J22").Value = 899.5
$ws.Range("K22").Value = 501
$ws.Range("L22").Value = 899.5
$ws.Range("M22").Value = -328
$ws.Range("N22").Value = -1245.5
$ws.Range("H107").Value = 9549.1875
$ws.Range("I107").Value = 11313.875
$ws.Range("J107").Value = 7784.5
$ws.Range("K107").Value = 11313.875
$ws.Range("L107").Value = 7784.5
$ws.Range("M107").Value = -9393.875
$ws.Range("N107").Value = -11624.5
$ws.Range("H128").Value = 3249.5
$ws.Range("I128").Value = 3249.5
$ws.Range("K128").Value = 9748.5
$ws.Range("M128").Value = -7258.5
$ws.Range("H135").Value = 39769.23
$ws.Range("J135").Value = 39769.23
$ws.Range("L135").Value = 39769.23
$ws.Range("N135").Value = -49909.23

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 3495
$ws.Range("I31").Value = 2993.3333
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 2993.3333
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -2698.3333
$ws.Range("N31").Value = -5590
$ws.Range("H34").Value = 3495
$ws.Range("I34").Value = 2993.3333
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 2993.3333
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -2791.3333
$ws.Range("N34").Value = -5404
$ws.Range("H57").Value = 44899
$ws.Range("J57").Value = 44899
$ws.Range("L57").Value = 44899
$ws.Range("N57").Value = -46019
$ws.Range("H92").Value = 19249.5
$ws.Range("J92").Value = 19249.5
$ws.Range("L92").Value = 19249.5
$ws.Range("N92").Value = -24241.5
$ws.Range("H95").Value = 15000
$ws.Range("J95").Value = 15000
$ws.Range("L95").Value = 15000
$ws.Range("N95").Value = -20492
$ws.Range("H96").Value = 10383.571
$ws.Range("J96").Value = 10383.571
$ws.Range("L96").Value = 10383.571
$ws.Range("N96").Value = -15875.571
$ws.Range("H99").Value = 1554
$ws.Range("I99").Value = 1554
$ws.Range("K99").Value = 1554
$ws.Range("M99").Value = -56
$ws.Range("H122").Value = 3971.6428
$ws.Range("I122").Value = 2445.6875
$ws.Range("K122").Value = 7337.0625
$ws.Range("M122").Value = -4887.0625
$ws.Range("H126").Value = 1554
$ws.Range("I126").Value = 1554
$ws.Range("K126").Value = 4662
$ws.Range("M126").Value = -2192

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H34").Value = 1514.3
$ws.Range("I34").Value = 160
$ws.Range("J34").Value = 2417.1667
$ws.Range("K34").Value = 480
$ws.Range("L34").Value = 7251.500100000001
$ws.Range("M34").Value = -396
$ws.Range("N34").Value = -7419.500100000001
$ws.Range("H47").Value = 1928.75
$ws.Range("I47").Value = 3257.5
$ws.Range("J47").Value = 1264.375
$ws.Range("K47").Value = 9772.5
$ws.Range("L47").Value = 3793.125
$ws.Range("M47").Value = -9341.5
$ws.Range("N47").Value = -4655.125
$ws.Range("H52").Value = 3500
$ws.Range("J52").Value = 3500
$ws.Range("L52").Value = 10500
$ws.Range("N52").Value = -11032
$ws.Range("H81").Value = 860
$ws.Range("I81").Value = 860
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2580
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1457
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 860
$ws.Range("I84").Value = 860
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 7740
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -2124
$ws.Range("N84").ClearContents()
$ws.Range("H99").Value = 32134.666
$ws.Range("I99").Value = 2020.875
$ws.Range("K99").Value = 6062.625
$ws.Range("M99").Value = -3816.625
$ws.Range("H119").Value = 4558
$ws.Range("I119").Value = 3860.2856
$ws.Range("K119").Value = 11580.8568
$ws.Range("M119").Value = -6742.856800000001
$ws.Range("H122").Value = 7709.5454
$ws.Range("I122").Value = 316.66666
$ws.Range("J122").Value = 10481.875
$ws.Range("K122").Value = 2849.99994
$ws.Range("L122").Value = 94336.875
$ws.Range("M122").Value = -399.9999399999997
$ws.Range("N122").Value = -99236.875
$ws.Range("H132").Value = 2009.1
$ws.Range("I132").Value = 1970.2858
$ws.Range("J132").Value = 2099.6667
$ws.Range("K132").Value = 17732.5722
$ws.Range("L132").Value = 18897.0003
$ws.Range("M132").Value = -15202.5722
$ws.Range("N132").Value = -23957.0003
$ws.Range("H137").Value = 2033.1666
$ws.Range("I137").Value = 1763.4546
$ws.Range("K137").Value = 5290.3638
$ws.Range("M137").Value = -190.3638000000001
$ws.Range("H140").Value = 2704.348
$ws.Range("I140").Value = 2478.889
$ws.Range("J140").Value = 2849.2856
$ws.Range("K140").Value = 7436.667
$ws.Range("L140").Value = 8547.856800000001
$ws.Range("M140").Value = -2256.667
$ws.Range("N140").Value = -18907.8568

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H43").Value = 20250
$ws.Range("I43").Value = 20250
$ws.Range("K43").Value = 20250
$ws.Range("M43").Value = -20099
$ws.Range("H80").Value = 54590.816
$ws.Range("I80").Value = 141701.38
$ws.Range("J80").Value = 4813.357
$ws.Range("K80").Value = 141701.38
$ws.Range("L80").Value = 4813.357
$ws.Range("M80").Value = -140703.38
$ws.Range("N80").Value = -6809.357
$ws.Range("H83").Value = 54590.816
$ws.Range("I83").Value = 141701.38
$ws.Range("J83").Value = 4813.357
$ws.Range("K83").Value = 708506.9
$ws.Range("L83").Value = 24066.785
$ws.Range("M83").Value = -703514.9
$ws.Range("N83").Value = -34050.785
$ws.Range("H88").Value = 29909.092
$ws.Range("J88").Value = 29909.092
$ws.Range("L88").Value = 29909.092
$ws.Range("N88").Value = -30811.092
$ws.Range("H91").Value = 29909.092
$ws.Range("J91").Value = 29909.092
$ws.Range("L91").Value = 29909.092
$ws.Range("N91").Value = -33029.092

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 28439.25
$ws.Range("I40").Value = 3930.5715
$ws.Range("K40").Value = 3930.5715
$ws.Range("M40").Value = -3794.5715
$ws.Range("H46").Value = 2329.8
$ws.Range("I46").Value = 1343.8
$ws.Range("K46").Value = 1343.8
$ws.Range("M46").Value = -1155.8
$ws.Range("H122").Value = 4587
$ws.Range("I122").Value = 1976.6666
$ws.Range("K122").Value = 5929.9998
$ws.Range("M122").Value = -3479.9998

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H56").Value = 29500
$ws.Range("J56").Value = 29500
$ws.Range("L56").Value = 29500
$ws.Range("N56").Value = -30928
$ws.Range("H94").Value = 20165
$ws.Range("J94").Value = 20165
$ws.Range("L94").Value = 20165
$ws.Range("N94").Value = -21967
$ws.Range("H107").Value = 608
$ws.Range("I107").Value = 612
$ws.Range("K107").Value = 1836
$ws.Range("M107").Value = 84
$ws.Range("H109").Value = 18578.947
$ws.Range("J109").Value = 18578.947
$ws.Range("L109").Value = 18578.947
$ws.Range("N109").Value = -21352.947
$ws.Range("H122").Value = 2206.5264
$ws.Range("I122").Value = 2120.5625
$ws.Range("K122").Value = 6361.6875
$ws.Range("M122").Value = -3911.6875
$ws.Range("H125").Value = 41181.816
$ws.Range("J125").Value = 41181.816
$ws.Range("L125").Value = 41181.816
$ws.Range("N125").Value = -51021.816

